$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 120 (pushes old rows 120..161 down to 122..163,
# and copies formatting -- including the date-number-format style -- from the row above).
$ws.Rows("120:121").Insert()

# New row 120: same as what later became row 122's "old" neighbour, with a new
# sample date (44559) and a revised volume.
$ws.Range("A120").Value = 8
$ws.Range("B120").Value = "Terminal La Palmera de La Serena"
$ws.Range("C120").Value = "Coquimbo"
$ws.Range("D120").Value = 44559
$ws.Range("E120").Value = 4
$ws.Range("F120").Value = 100112031
$ws.Range("G120").Value = "Poroto verde"
$ws.Range("H120").Value = "Magnum"
$ws.Range("I120").Value = "Primera"
$ws.Range("J120").Value = 500
$ws.Range("K120").Value = 24000
$ws.Range("L120").Value = 25000
$ws.Range("M120").Value = 24500
$ws.Range("N120").Value = "$/malla 25 kilos"
$ws.Range("O120").Value = "Provincia de Limarí"
$ws.Range("P120").Value = 980
$ws.Range("Q120").Value = 25
$ws.Range("R120").Value = "Hortaliza"

# New row 121: second newly-inserted sample row.
$ws.Range("A121").Value = 8
$ws.Range("B121").Value = "Terminal La Palmera de La Serena"
$ws.Range("C121").Value = "Coquimbo"
$ws.Range("D121").Value = 44559
$ws.Range("E121").Value = 4
$ws.Range("F121").Value = 100112031
$ws.Range("G121").Value = "Poroto verde"
$ws.Range("H121").Value = "Sin especificar"
$ws.Range("I121").Value = "Primera"
$ws.Range("J121").Value = 400
$ws.Range("K121").Value = 28000
$ws.Range("L121").Value = 29000
$ws.Range("M121").Value = 28500
$ws.Range("N121").Value = "$/malla 25 kilos"
$ws.Range("O121").Value = "Provincia del Elquí"
$ws.Range("P121").Value = 1140
$ws.Range("Q121").Value = 25
$ws.Range("R121").Value = "Hortaliza"

# Make sure the date cells keep the date/time number format used throughout
# column D (style index "s=2" -> numFmt "YYYY-MM-DD HH:MM:SS").
$ws.Range("D120:D121").NumberFormat = "YYYY-MM-DD HH:MM:SS"
